# "Generate Report for Handoff" — refresh the Latest HO Xliff Generate Date /
# Latest Handoff Datetime timestamps for the d9f50b50 file (row 7) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-19 08:40:04"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-19 08:39:57"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-19 08:40:04"
